# Updates the date line and refreshes the multiplication problems
# with a new set of two-digit x two-digit factors.
$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("2026-01-15 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-16 Friday", 2)
$null = $d.Content.Find.Execute("50×93=", $true, $false, $false, $false, $false, $true, 1, $false, "53×80=", 2)
$null = $d.Content.Find.Execute("56×56=", $true, $false, $false, $false, $false, $true, 1, $false, "78×90=", 2)
$null = $d.Content.Find.Execute("66×69=", $true, $false, $false, $false, $false, $true, 1, $false, "87×88=", 2)
$null = $d.Content.Find.Execute("27×59=", $true, $false, $false, $false, $false, $true, 1, $false, "51×93=", 2)
$null = $d.Content.Find.Execute("49×33=", $true, $false, $false, $false, $false, $true, 1, $false, "16×78=", 2)
$null = $d.Content.Find.Execute("67×17=", $true, $false, $false, $false, $false, $true, 1, $false, "53×53=", 2)
$null = $d.Content.Find.Execute("41×13=", $true, $false, $false, $false, $false, $true, 1, $false, "82×21=", 2)
$null = $d.Content.Find.Execute("27×17=", $true, $false, $false, $false, $false, $true, 1, $false, "98×48=", 2)
$null = $d.Content.Find.Execute("25×83=", $true, $false, $false, $false, $false, $true, 1, $false, "57×82=", 2)
$null = $d.Content.Find.Execute("31×61=", $true, $false, $false, $false, $false, $true, 1, $false, "36×26=", 2)
$null = $d.Content.Find.Execute("87×43=", $true, $false, $false, $false, $false, $true, 1, $false, "67×37=", 2)
$null = $d.Content.Find.Execute("48×74=", $true, $false, $false, $false, $false, $true, 1, $false, "52×34=", 2)
$null = $d.Content.Find.Execute("20×48=", $true, $false, $false, $false, $false, $true, 1, $false, "53×35=", 2)
$null = $d.Content.Find.Execute("85×21=", $true, $false, $false, $false, $false, $true, 1, $false, "50×61=", 2)
$null = $d.Content.Find.Execute("79×29=", $true, $false, $false, $false, $false, $true, 1, $false, "25×35=", 2)
$null = $d.Content.Find.Execute("11×34=", $true, $false, $false, $false, $false, $true, 1, $false, "15×22=", 2)
$null = $d.Content.Find.Execute("39×57=", $true, $false, $false, $false, $false, $true, 1, $false, "38×13=", 2)
$null = $d.Content.Find.Execute("69×47=", $true, $false, $false, $false, $false, $true, 1, $false, "27×40=", 2)
$null = $d.Content.Find.Execute("77×62=", $true, $false, $false, $false, $false, $true, 1, $false, "69×47=", 2)
$null = $d.Content.Find.Execute("37×76=", $true, $false, $false, $false, $false, $true, 1, $false, "46×41=", 2)
$null = $d.Content.Find.Execute("95×94=", $true, $false, $false, $false, $false, $true, 1, $false, "40×95=", 2)
$null = $d.Content.Find.Execute("73×66=", $true, $false, $false, $false, $false, $true, 1, $false, "55×26=", 2)
$null = $d.Content.Find.Execute("15×18=", $true, $false, $false, $false, $false, $true, 1, $false, "62×31=", 2)
$null = $d.Content.Find.Execute("66×85=", $true, $false, $false, $false, $false, $true, 1, $false, "42×26=", 2)
$null = $d.Content.Find.Execute("23×77=", $true, $false, $false, $false, $false, $true, 1, $false, "41×54=", 2)
